# Add a reviewer comment "Add value" anchored on the "YYY" placeholder
# inside " (YYY)", matching the target diff:
#   " (YYY)" -> " (" + [commentRangeStart]YYY[commentRangeEnd][commentRef] + ")"
# Word's COM model renumbers bookmark/comment w:id values and regenerates
# commentsExtended/commentsExtensible/commentsIds automatically when a
# comment is added via Comments.Add, so a single Add call keeps every part
# of the package in sync.

$word.UserName = "Federico Tartarini"
$word.UserInitials = "f"

$d = $word.ActiveDocument

# Locate the "YYY" placeholder text (only occurrence in the document) and
# narrow the range down to exactly those three characters.
$target = $d.Content
$found = $target.Find.Execute("YYY", $false, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find placeholder 'YYY' to comment on"
}

# Anchor a new comment on the found range ("YYY"); Word wraps it with
# commentRangeStart/commentRangeEnd/commentReference automatically.
$comment = $d.Comments.Add($target, "Add value")
$comment.Initial = "f"
